$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels: eta_min -> Y_min, eta_max -> Y_max, eta -> Y
$ws.Range("D1").Value = "Y_min"
$ws.Range("E1").Value = "Y_max"
$ws.Range("F1").Value = "Y"

# Update the "obs" column values for all data rows: Z_rap -> Z_xsec_norm
$ws.Range("M2:M29").Value = "Z_xsec_norm"

# Widen column M (13) to fit the longer label, matching the new column layout
$ws.Columns.Item(13).ColumnWidth = 14.5

# Update the active selection to reflect the edited area
$ws.Range("M31").Select()
